$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 467 (existing rows 467:484 shift down to 469:486)
$ws.Rows("467:468").Insert()

# New row 467 data
$ws.Cells.Item(467,1).Value  = 4
$ws.Cells.Item(467,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(467,3).Value  = "Los Lagos"
$ws.Cells.Item(467,4).Value  = 45041
$ws.Cells.Item(467,5).Value  = 10
$ws.Cells.Item(467,6).Value  = "Fruta"
$ws.Cells.Item(467,7).Value  = 100104
$ws.Cells.Item(467,8).Value  = "Frutos de pepita"
$ws.Cells.Item(467,9).Value  = 100104005
$ws.Cells.Item(467,10).Value = "Pera"
$ws.Cells.Item(467,11).Value = "Packham's Triumph"
$ws.Cells.Item(467,12).Value = "Primera"
$ws.Cells.Item(467,13).Value = 400
$ws.Cells.Item(467,14).Value = 18000
$ws.Cells.Item(467,15).Value = 19000
$ws.Cells.Item(467,16).Value = 18500
$ws.Cells.Item(467,17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(467,18).Value = "Región de O'Higgins"
$ws.Cells.Item(467,19).Value = 1233
$ws.Cells.Item(467,20).Value = 15

# New row 468 data
$ws.Cells.Item(468,1).Value  = 4
$ws.Cells.Item(468,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(468,3).Value  = "Los Lagos"
$ws.Cells.Item(468,4).Value  = 45041
$ws.Cells.Item(468,5).Value  = 10
$ws.Cells.Item(468,6).Value  = "Fruta"
$ws.Cells.Item(468,7).Value  = 100104
$ws.Cells.Item(468,8).Value  = "Frutos de pepita"
$ws.Cells.Item(468,9).Value  = 100104005
$ws.Cells.Item(468,10).Value = "Pera"
$ws.Cells.Item(468,11).Value = "Packham's Triumph"
$ws.Cells.Item(468,12).Value = "Segunda"
$ws.Cells.Item(468,13).Value = 200
$ws.Cells.Item(468,14).Value = 16000
$ws.Cells.Item(468,15).Value = 16000
$ws.Cells.Item(468,16).Value = 16000
$ws.Cells.Item(468,17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(468,18).Value = "Región de O'Higgins"
$ws.Cells.Item(468,19).Value = 1067
$ws.Cells.Item(468,20).Value = 15
